# Fix issue with "Trenching" input to installation module.
# Replace incorrectly removed "ROV class" input.
# Add DateTimeDict for date outputs from installation module.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rov")

# --- Insert a new "Class" column before the old column C (Depth Rating) ---
# This shifts every existing column from C..P one place to the right (D..Q),
# preserving all existing values/styles, and re-sizes the used dimension.
$ws.Columns("C:C").Insert()

# New column header
$ws.Range("C1").Value = "Class"

# New column values: the previously-generic "Name" values ("Inspection class" /
# "Workclass") move here, while column B gets unique per-row names.
$ws.Range("B2").Value = "Inspection class 1"
$ws.Range("C2").Value = "Inspection class"

$ws.Range("B3").Value = "Inspection class 2"
$ws.Range("C3").Value = "Inspection class"

$ws.Range("B4").Value = "Workclass 1"
$ws.Range("C4").Value = "Workclass"

$ws.Range("B5").Value = "Workclass 2"
$ws.Range("C5").Value = "Workclass"

$ws.Range("B6").Value = "Workclass 3"
$ws.Range("C6").Value = "Workclass"

# Restore the new column's width/formatting (Excel's Insert copies the width
# of the column to its left by default; set it to the width used in the
# saved workbook).
$ws.Range("C1").ColumnWidth = 18.81640625

# --- Fix up the cell comments, which Excel does not automatically re-anchor
# when a column is inserted (the underlying cells moved, but the comment
# anchors did not) ---
$oldRefs = @("N2", "P2", "P3", "G5", "N5", "O5", "P5")
foreach ($r in $oldRefs) {
    $cmt = $ws.Range($r).Comment
    if ($cmt) { $cmt.Delete() }
}

$nl = [char]10

$c = $ws.Range("O2").AddComment("Author:" + $nl + "cost of a different oc class rov")
$c = $ws.Range("Q2").AddComment("Author:" + $nl + "labour ratesfrom oceaneering in EUR 2015")
$c = $ws.Range("Q3").AddComment("Author:" + $nl + "One technicians @ Euros 1800 each")
$c = $ws.Range("H5").AddComment("Author:" + $nl + "Including the payload of 400kg")
$c = $ws.Range("O5").AddComment("Author:" + $nl + "Original value from fugro probably too high compared to market rates. Discussions with  4C Offshore's subsea consultant Dr. Ron Haynes (who has a lot of experience managing subsea cable contracts) suggested " + [char]0x00A3 + "GBP11K more appropriate market rate for the spread, with reputable teams such as Fugro charging a premium")
$c = $ws.Range("P5").AddComment("Author:" + $nl + "1 supervisor")
$c = $ws.Range("Q5").AddComment("Author:" + $nl + "4 technicians at EUR 1100 each")

# --- Selection / active-cell bookkeeping on the "rov" sheet ---
$ws.Range("B8").Select()

# --- Make "rov" (first tab) the active sheet instead of "cable_burial" ---
$ws.Activate()
